$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row additions (columns AD, AE, AF) - match style of existing header cell (bold, border, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-45: team record values (constant across all rows)
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 88   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
